# Applies the diff described for Archivos/HorasDedicadas.xlsx:
#  - Adds 3 new rows (43, 44, 45) of work-log data to the sheet / Tabla1 table.
#  - Extends Tabla1 (ListObject) range from D4:F42 to D4:F45.
#  - Updates the summary formulas in I5 / J5 (recalculated automatically).
#  - Updates the sheet selection / scroll position to match the new last cell (F45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1. Write the new data values FIRST (before touching any formatting/formula),
#    so the calculation engine picks up the new cell contents correctly.
$ws.Range("D43").Value = 44276
$ws.Range("E43").Value = 75
$ws.Range("F43").Value = "Crear listView de la pokedex y actividad intermedia."

$ws.Range("D44").Value = 44276
$ws.Range("E44").Value = 40
$ws.Range("F44").Value = "Cambiar metodos de conexión para poder unirse en mitad de la partida y abandonar."

$ws.Range("D45").Value = 44276
$ws.Range("E45").Value = 95
$ws.Range("F45").Value = "Cambiar conexiones para poder unirse y salirse en cualquier momento y evitar fallos."

# 2. Re-assert the summary formulas so they recalculate against the new data.
$ws.Range("I5").Formula = "=SUM(E5:E45)"
$ws.Range("J5").Formula = "=CONVERT(I5,""mn"",""hr"")"

# 3. Copy the formatting (number format / borders / alignment) from the last
#    existing data row (42) onto the three new rows.
$ws.Range("D42:F42").Copy()
$ws.Range("D43:F43").PasteSpecial(-4122)
$ws.Range("D42:F42").Copy()
$ws.Range("D44:F44").PasteSpecial(-4122)
$ws.Range("D42:F42").Copy()
$ws.Range("D45:F45").PasteSpecial(-4122)

# 4. Grow the table ("Tabla1") so it covers the new rows as well.
$lo.Resize($ws.Range("D4:F45"))

# 5. Update the view: select F45 (removes the stale topLeftCell scroll anchor
#    and matches the author's final cursor position).
$ws.Range("F45").Select()
